$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFCB109")

# -----------------------------------------------------------------
# Row 3 - fill in new data for file 20160718_213034
# -----------------------------------------------------------------
# Numeric / already-known-string cells first (order does not affect
# the shared-string table since they are numbers or reuse an
# existing shared string).
$ws.Cells.Item(3, 1).Value = 20160718
$ws.Cells.Item(3, 2).Value = 213034
$ws.Cells.Item(3, 3).Value = 0.55
$ws.Cells.Item(3, 4).Value = "NA"
$ws.Cells.Item(3, 5).Value = 1.3986
$ws.Cells.Item(3, 6).Value = 0.077093
$ws.Cells.Item(3, 9).Value = 0.5

# New-string cells: write them in the same order they were first
# introduced by the author so the shared-string table indices line
# up with the target workbook.
$ws.Cells.Item(3, 16).Value = "very skewed"
$ws.Cells.Item(3, 8).Value = "1.0-2.0"
$ws.Cells.Item(3, 8).NumberFormat = "d-mmm"
$ws.Cells.Item(3, 7).Value = "0.05-0.15"
$ws.Cells.Item(3, 10).Value = "na"
$ws.Cells.Item(3, 11).Value = "na"

# -----------------------------------------------------------------
# Row 4 - new row for file 20160721_155707
# -----------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = 20160721
$ws.Cells.Item(4, 2).Value = 155707
$ws.Cells.Item(4, 3).Value = 0.55
$ws.Cells.Item(4, 5).Value = 1.009
$ws.Cells.Item(4, 6).Value = 0.058626
$ws.Cells.Item(4, 9).Value = 0.5
$ws.Cells.Item(4, 10).Value = 1.5672
$ws.Cells.Item(4, 11).Value = 0.090564
$ws.Cells.Item(4, 12).Value = 0.24626
$ws.Cells.Item(4, 13).Value = 0.11241
$ws.Cells.Item(4, 14).Value = 0.18114
$ws.Cells.Item(4, 15).Value = 0.10228

$ws.Cells.Item(4, 16).Value = "moving laser a lot, is this acceptable alignment?"
$ws.Cells.Item(4, 17).Value = "select out singlet beads, a few 20um and doublets in there"
$ws.Cells.Item(4, 8).Value = "0.75-1.15"
$ws.Cells.Item(4, 7).Value = "0.045-0.07"
$ws.Cells.Item(4, 4).Value = "na"

# -----------------------------------------------------------------
# Update the active cell selection to O6, matching the author's
# final cursor position in the saved workbook.
# -----------------------------------------------------------------
$ws.Range("O6").Select() | Out-Null
